$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data1 = New-Object 'object[,]' 24,5
$data1[0,0] = 1.02
$data1[0,1] = 1.025648509051531
$data1[0,2] = 1.031663787263372
$data1[0,3] = 0.9926147277508489
$data1[0,4] = 1.037012208397243
$data1[1,0] = 1.02
$data1[1,1] = 1.026413331397954
$data1[1,2] = 1.032245844570696
$data1[1,3] = 0.9936372048519304
$data1[1,4] = 1.038058878356229
$data1[2,0] = 1.02
$data1[2,1] = 1.02690849058097
$data1[2,2] = 1.032622597489315
$data1[2,3] = 0.9942998659930995
$data1[2,4] = 1.038736779765034
$data1[3,0] = 1.02
$data1[3,1] = 1.027116718119278
$data1[3,2] = 1.032781012211452
$data1[3,3] = 0.9945786998346017
$data1[3,4] = 1.03902192076283
$data1[4,0] = 1.02
$data1[4,1] = 1.027151684066425
$data1[4,2] = 1.032807612322736
$data1[4,3] = 0.9946255319796338
$data1[4,4] = 1.039069806031055
$data1[5,0] = 1.02
$data1[5,1] = 1.02691127268341
$data1[5,2] = 1.032624714127427
$data1[5,3] = 0.9943035907982488
$data1[5,4] = 1.038740589239962
$data1[6,0] = 1.02
$data1[6,1] = 1.025906927932406
$data1[6,2] = 1.031860469824288
$data1[6,3] = 0.9929600610674301
$data1[6,4] = 1.037365803126916
$data1[7,0] = 1.02
$data1[7,1] = 1.024139264464532
$data1[7,2] = 1.030514787082385
$data1[7,3] = 0.9906006454969559
$data1[7,4] = 1.034948169238433
$data1[8,0] = 1.02
$data1[8,1] = 1.022962338711782
$data1[8,2] = 1.029618443367496
$data1[8,3] = 0.989033133672735
$data1[8,4] = 1.033339775844418
$data1[9,0] = 1.02
$data1[9,1] = 1.022453095595742
$data1[9,2] = 1.029230521553725
$data1[9,3] = 0.988355674866747
$data1[9,4] = 1.032644133013884
$data1[10,0] = 1.02
$data1[10,1] = 1.022263997655393
$data1[10,2] = 1.029086461943331
$data1[10,3] = 0.9881042295826724
$data1[10,4] = 1.032385861685054
$data1[11,0] = 1.02
$data1[11,1] = 1.022304557173483
$data1[11,2] = 1.029117361752489
$data1[11,3] = 0.9881581567098651
$data1[11,4] = 1.03244125625437
$data1[12,0] = 1.02
$data1[12,1] = 1.022437463516683
$data1[12,2] = 1.029218612882758
$data1[12,3] = 0.9883348863814464
$data1[12,4] = 1.032622781743801
$data1[13,0] = 1.02
$data1[13,1] = 1.022519359182567
$data1[13,2] = 1.029281001311592
$data1[13,3] = 0.9884438009545853
$data1[13,4] = 1.032734641701
$data1[14,0] = 1.02
$data1[14,1] = 1.022996143498615
$data1[14,2] = 1.02964419283381
$data1[14,3] = 0.9890781214508737
$data1[14,4] = 1.033385960343042
$data1[15,0] = 1.02
$data1[15,1] = 1.023295318768535
$data1[15,2] = 1.02987206823144
$data1[15,3] = 0.989476357848556
$data1[15,4] = 1.03379473069221
$data1[16,0] = 1.02
$data1[16,1] = 1.023469858577912
$data1[16,2] = 1.030005003372908
$data1[16,3] = 0.9897087662937556
$data1[16,4] = 1.034033236854907
$data1[17,0] = 1.02
$data1[17,1] = 1.023529378191747
$data1[17,2] = 1.030050334040723
$data1[17,3] = 0.9897880325774034
$data1[17,4] = 1.034114574383222
$data1[18,0] = 1.02
$data1[18,1] = 1.023263216365469
$data1[18,2] = 1.029847617352096
$data1[18,3] = 0.9894336180360679
$data1[18,4] = 1.033750865496877
$data1[19,0] = 1.02
$data1[19,1] = 1.022398324308282
$data1[19,2] = 1.029188796056915
$data1[19,3] = 0.9882828385668249
$data1[19,4] = 1.032569323653453
$data1[20,0] = 1.02
$data1[20,1] = 1.021854866504995
$data1[20,2] = 1.028774753887017
$data1[20,3] = 0.9875604150241495
$data1[20,4] = 1.031827144177039
$data1[21,0] = 1.02
$data1[21,1] = 1.022142931589543
$data1[21,2] = 1.028994227401309
$data1[21,3] = 0.9879432794643023
$data1[21,4] = 1.032220520649421
$data1[22,0] = 1.02
$data1[22,1] = 1.023277721948717
$data1[22,2] = 1.029858665592632
$data1[22,3] = 0.9894529299347244
$data1[22,4] = 1.033770686051573
$data1[23,0] = 1.02
$data1[23,1] = 1.024595987196539
$data1[23,2] = 1.030862548246915
$data1[23,3] = 0.9912096547607049
$data1[23,4] = 1.03557259714058
$ws.Range("B2:F25").Value = $data1

$data2 = New-Object 'object[,]' 24,6
$data2[0,0] = 1.035346229762267
$data2[0,1] = 1.030816355624612
$data2[0,2] = 1.034471253199687
$data2[0,3] = 0.9955398523336033
$data2[0,4] = 1.039804308982901
$data2[0,5] = 1.014326721090837
$data2[1,0] = 1.035540994721242
$data2[1,1] = 1.031221554070435
$data2[1,2] = 1.034862528631744
$data2[1,3] = 0.9963617723202692
$data2[1,4] = 1.040660071045232
$data2[1,5] = 1.014460439170337
$data2[2,0] = 1.035665774165116
$data2[2,1] = 1.031483335274174
$data2[2,2] = 1.035115133404566
$data2[2,3] = 0.9968940712668345
$data2[2,4] = 1.041213843072462
$data2[2,5] = 1.014546821393938
$data2[3,0] = 1.035717932310691
$data2[3,1] = 1.031593289070143
$data2[3,2] = 1.035221189288547
$data2[3,3] = 0.997117960005301
$data2[3,4] = 1.041446656644217
$data2[3,5] = 1.014583102049638
$data2[4,0] = 1.035726672330742
$data2[4,1] = 1.031611744959388
$data2[4,2] = 1.035238988344849
$data2[4,3] = 0.9971555583673453
$data2[4,4] = 1.041485747496796
$data2[4,5] = 1.014589191701821
$data2[5,0] = 1.035666472281091
$data2[5,1] = 1.03148480487257
$data2[5,2] = 1.035116551076705
$data2[5,3] = 0.9968970624462087
$data2[5,4] = 1.04121695390934
$data2[5,5] = 1.014547306313579
$data2[6,0] = 1.035412309159276
$data2[6,1] = 1.030953378458003
$data2[6,2] = 1.034603605276426
$data2[6,3] = 0.995817528259106
$data2[6,4] = 1.04009350957546
$data2[6,5] = 1.014371940992642
$data2[7,0] = 1.034954926092173
$data2[7,1] = 1.030013848355922
$data2[7,2] = 1.033695362772354
$data2[7,3] = 0.9939188001724441
$data2[7,4] = 1.038114182045377
$data2[7,5] = 1.014061850817192
$data2[8,0] = 1.034643648885162
$data2[8,1] = 1.029385479550304
$data2[8,2] = 1.03308699851809
$data2[8,3] = 0.9926553831429383
$data2[8,4] = 1.036794906089588
$data2[8,5] = 1.013854423200603
$data2[9,0] = 1.034507363871323
$data2[9,1] = 1.029112923722397
$data2[9,2] = 1.0328229034554
$data2[9,3] = 0.9921088820399291
$data2[9,4] = 1.036223722073807
$data2[9,5] = 1.013764442983252
$data2[10,0] = 1.034456516884185
$data2[10,1] = 1.02901161505776
$data2[10,2] = 1.032724707335139
$data2[10,3] = 0.9919059725120875
$data2[10,4] = 1.036011570472438
$data2[10,5] = 1.013730996235242
$data2[11,0] = 1.034467433891133
$data2[11,1] = 1.029033349229667
$data2[11,2] = 1.03274577521913
$data2[11,3] = 0.9919494934313052
$data2[11,4] = 1.036057077173183
$data2[11,5] = 1.013738171760634
$data2[12,0] = 1.034503165422653
$data2[12,1] = 1.029104550924297
$data2[12,2] = 1.032814788555948
$data2[12,3] = 0.9920921077337197
$data2[12,4] = 1.036206185315047
$data2[12,5] = 1.013761678753992
$data2[13,0] = 1.03452515104865
$data2[13,1] = 1.029148411486449
$data2[13,2] = 1.032857296804687
$data2[13,3] = 0.9921799884222134
$data2[13,4] = 1.036298057322499
$data2[13,5] = 1.013776159006975
$data2[14,0] = 1.034652662129156
$data2[14,1] = 1.029403558409689
$data2[14,2] = 1.03310451162994
$data2[14,3] = 0.9926916645766087
$data2[14,4] = 1.036832815266291
$data2[14,5] = 1.013860391493774
$data2[15,0] = 1.034732245369298
$data2[15,1] = 1.029563480757251
$data2[15,2] = 1.033259404540634
$data2[15,3] = 0.9930127773699352
$data2[15,4] = 1.037168274473496
$data2[15,5] = 1.013913185000255
$data2[16,0] = 1.034778520159653
$data2[16,1] = 1.029656715534816
$data2[16,2] = 1.033349686300734
$data2[16,3] = 0.9932001317071769
$data2[16,4] = 1.037363948962091
$data2[16,5] = 1.013943962838448
$data2[17,0] = 1.034794274091528
$data2[17,1] = 1.029688498483194
$data2[17,2] = 1.033380459047769
$data2[17,3] = 0.9932640239640975
$data2[17,4] = 1.037430670034133
$data2[17,5] = 1.01395445460859
$data2[18,0] = 1.034723721809813
$data2[18,1] = 1.02954632727304
$data2[18,2] = 1.033242792676711
$data2[18,3] = 0.9929783193494215
$data2[18,4] = 1.037132282164672
$data2[18,5] = 1.013907522380177
$data2[19,0] = 1.034492649573574
$data2[19,1] = 1.029083585704608
$data2[19,2] = 1.032794468584965
$data2[19,3] = 0.9920501090198102
$data2[19,4] = 1.036162276357106
$data2[19,5] = 1.013754757192481
$data2[20,0] = 1.03434606578355
$data2[20,1] = 1.028792241298664
$data2[20,2] = 1.032512014562444
$data2[20,3] = 0.9914670000341481
$data2[20,4] = 1.035552462608399
$data2[20,5] = 1.013658568459161
$data2[21,0] = 1.034423895574108
$data2[21,1] = 1.028946726096821
$data2[21,2] = 1.032661802907312
$data2[21,3] = 0.991776070289318
$data2[21,4] = 1.035875729731026
$data2[21,5] = 1.013709573003027
$data2[22,0] = 1.034727573686828
$data2[22,1] = 1.029554078334034
$data2[22,2] = 1.033250299062865
$data2[22,3] = 0.9929938892766442
$data2[22,4] = 1.03714854551876
$data2[22,5] = 1.013910081123061
$data2[23,0] = 1.035074293816397
$data2[23,1] = 1.030257099123394
$data2[23,2] = 1.033930675816643
$data2[23,3] = 0.9944092447426414
$data2[23,4] = 1.038625841258271
$data2[23,5] = 1.014142141644606
$ws.Range("I2:N25").Value = $data2

Write-Output "Updated vm_pu values for case with 380 kV"
